$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "NG Price" -> "NG_Price"
$ws.Range("B1").Value = "NG_Price"

# The header cell A1 picks up a plain numeric format (matches the
# column-wide formatting cleanup done on column A). Applied first so it
# lands on its own style slot.
$ws.Cells.Item(1, 1).NumberFormat = "0.00"

# Column A currently holds dates (1997-06-30 .. 2019-06-30) displayed as
# plain years via a custom "yyyy" number format. Replace the values with
# plain integer years and format them with a simple integer ("0") number
# format instead of the date-based one.
$years = 1997..2019
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "0"
}

# Selection moved to E7 in the saved file
$ws.Range("E7").Select()
